# Updated template and flag map
# Insert two new rows for "administration_term" / "administration_term_units"
# right after the existing "administration_route" row (row 16), shifting the
# rest of the "studies" / "subjects" / "series" / "conc_time_values" block
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 17:18 (everything from row 17 down shifts by 2)
$ws.Rows("17:18").Insert()

# Populate the two new rows (row 18 first so the shared-string table gets
# the same allocation order as the authored workbook: administration_term_units
# before administration_term)
$ws.Range("A18").Value2 = "studies"
$ws.Range("B18").Value2 = "administration_term_units"
$ws.Range("C18").Value2 = "administration_term_units"

$ws.Range("A17").Value2 = "studies"
$ws.Range("B17").Value2 = "administration_term"
$ws.Range("C17").Value2 = "administration_term"

# Re-establish the autofilter over the now-larger range (A1:C78).
# Toggling AutoFilter() with no args flips AutoFilterMode, so make sure we
# end up re-applying it cleanly regardless of its current on/off state.
if ($ws.AutoFilterMode) { $ws.Range("A1:C78").AutoFilter() }
$ws.Range("A1:C78").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# resized autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$78"
    }
}

# Update the active selection to B2 (also clears the stale topLeftCell scroll
# position left over from before the edit).
$ws.Range("B2").Select()
